$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JAN-2021")

# Copy formatting from the last existing data row (row 6) down into the
# three new rows so the new cells pick up the same styles (borders,
# alignment, number format, wrap text) without creating new style records.
$ws.Range("A6:G6").Copy()
$ws.Range("A7:G9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 7 (Jan 11, 2021)
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 44207
$ws.Range("C7").Value = "Selenium log Files(Sony, Samsung)"
$ws.Range("D7").Value = "Selenium log file Testing (QMVAR TO GSPN)"
$ws.Range("G7").Value = "Unable to import, Issue Find(Sony Daily Claim)"

# Row 8 (Jan 12, 2021)
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 44208
$ws.Range("C8").Value = "Selenium log Files(Sony, Samsung), GitHub "
$ws.Range("G8").Value = "GitHub Pull and Push"

# Row 9 (Jan 13, 2021)
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 44209
$ws.Range("C9").Value = "Selenium log Files(Sony, Samsung), Git Hub"
$ws.Range("G9").Value = "Git Hub Admin Testing(User Privilage)"

# Status/Application column - filled in after the Task/Comments columns,
# matching the order the new shared strings were originally introduced.
$ws.Range("D8").Value = "Selenium log file Testing (QMVAR TO GSPN), Git Hub"
$ws.Range("D9").Value = "Selenium log file Testing (QMVAR TO GSPN), Git Hub"

# % of completion / Status columns
$ws.Range("E7").Value = 1
$ws.Range("E8").Value = 1
$ws.Range("E9").Value = 1
$ws.Range("F7").Value = "Completed"
$ws.Range("F8").Value = "Completed"
$ws.Range("F9").Value = "Completed"

# Match the taller wrapped-text row height used by the other data rows.
$ws.Rows.Item(7).RowHeight = 28.8
$ws.Rows.Item(8).RowHeight = 28.8
$ws.Rows.Item(9).RowHeight = 28.8

# Move the active selection to F9, as left by the author after data entry.
$ws.Range("F9").Select()
